# Update Name of Algo
# Apply updated numeric values to specific cells on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3"  = 4.639200000000004
    "E3"  = 12.4138
    "C6"  = -11.7464
    "C7"  = -12.2839
    "B8"  = 5.221899999999996
    "C8"  = -11.1703
    "E8"  = 13.28969999999999
    "D11" = -8.645299999999992
    "A12" = -22.7471
    "B12" = 5.5827
    "B14" = 9.101000000000001
    "D14" = -8.764299999999995
    "C19" = -13.0293
    "D19" = -8.525099999999991
    "E19" = 13.3709
    "C21" = -13.02740000000001
    "D21" = -8.705999999999987
    "B22" = 4.878100000000005
    "C24" = -11.674
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
